$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.901.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.011.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.018.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

# Row 10
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.538.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.30%  "

# Row 14
$ws.Range("E14").Value = "  -3.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.000.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

# Row 16
$ws.Range("E16").Value = "  -0.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.013.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.60%  "

# Row 18
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "395.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.79%  "

# Row 23
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.465"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "

# Row 26
$ws.Range("E26").Value = "  -4.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0972"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "

# Row 29
$ws.Range("E29").Value = "  -0.64%  "

# Row 30
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("E31").Value = "  -0.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "

# Row 37
$ws.Range("E37").Value = "  +0.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.508.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.74%  "

# Row 39
$ws.Range("E39").Value = "  -2.01%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.666"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.99%  "

# Row 44
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0246"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.18%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0946"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.94%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.50"
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "264.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.16%  "

